$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of cell F8 (shared string used by F8) to the revised prompt text.
$newText = "设计任务/请设计10个任务并对每个任务提供要完成的主要步骤。  Final output are in the following format:     - answer 1     - answer 2     - answer 3"
$ws.Range("F8").Value2 = $newText

# Update the active cell / selection shown when the sheet was last saved.
$ws.Range("R4").Select()
